{"js": "// Add a new bulleted paragraph at the end of the document (same list as the\n// surrounding \"Prrafodelista\" bullets, but at the top ilvl=0 like the\n// section-header bullets), carrying the note about Jaime Marinas, and move\n// the \"_GoBack\" bookmark from the end of the previous paragraph onto the\n// newly inserted run.\n\nconst body = context.document.body;\n\n// The \"_GoBack\" bookmark currently sits (collapsed) at the end of the last\n// paragraph. Remove it so we can re-insert it around the new run further\n// down, matching Word's habit of re-anchoring _GoBack at the most recent\n// edit.\ncontext.document.deleteBookmark(\"_GoBack\");\n\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\n\nconst newParagraph = lastParagraph.insertParagraph(\n  \"Se considera la posibilidad de rescindir el contrato con Jaime Marinas debido a la nula aportaci\u00f3n y compromiso con el equipo. Se decide de enviar un email al responsable del proyecto, Andr\u00e9s Castillo, para ver las opciones que tiene el equipo de desarrollo para poder tomar una decisi\u00f3n.\",\n  \"After\"\n);\n\n// Same \"Prrafodelista\" list-paragraph style used throughout the document.\nnewParagraph.style = \"Prrafodelista\";\n\n// Join the existing list (numId=2) at the top level (ilvl=0), like the other\n// top-level bullets (e.g. \"Valoraci\u00f3n de la funcionalidad...\").\nnewParagraph.attachToList(2, 0);\n\n// Justify the paragraph text.\nnewParagraph.alignment = \"Justified\";\n\n// Re-insert the _GoBack bookmark wrapping the new run's text.\nnewParagraph.getRange(\"Content\").insertBookmark(\"_GoBack\");\n\nawait context.sync();\n", "ps1": "# Add a new bulleted paragraph at the end of the document (joining the same\n# list as the surrounding \"Prrafodelista\" bullets, but at the top level\n# ilvl=0 like the section-header bullets), carrying the note about Jaime\n# Marinas, and move the \"_GoBack\" bookmark from the end of the previous\n# paragraph onto the newly inserted run.\n\n$d = $word.ActiveDocument\n\n# The \"_GoBack\" bookmark currently sits (collapsed) at the end of the last\n# paragraph. Remove it so it can be re-inserted around the new run below,\n# matching Word's habit of re-anchoring _GoBack at the most recent edit.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n$lastParagraph = $d.Paragraphs.Last\n$r = $lastParagraph.Range\n$r.Collapse(0)  # wdCollapseEnd\n$r.InsertParagraphAfter()\n\n$newParagraph = $d.Paragraphs.Last\n$newParagraph.Range.Text = \"Se considera la posibilidad de rescindir el contrato con Jaime Marinas debido a la nula aportaci\u00f3n y compromiso con el equipo. Se decide de enviar un email al responsable del proyecto, Andr\u00e9s Castillo, para ver las opciones que tiene el equipo de desarrollo para poder tomar una decisi\u00f3n.\"\n\n# Join the existing list (numId=2) at the top level (ilvl=0, i.e.\n# ListLevelNumber 1 in the 1-based COM numbering), like the other top-level\n# bullets (e.g. \"Valoraci\u00f3n de la funcionalidad...\").\n$newParagraph.Range.ListFormat.ListLevelNumber = 1\n\n# Justify the paragraph text.\n$newParagraph.Alignment = 3  # wdAlignParagraphJustify\n\n# Re-insert the _GoBack bookmark wrapping the new run's text.\n$d.Bookmarks.Add(\"_GoBack\", $newParagraph.Range)\n"}
